$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7, shifting existing rows 7:49 down to 8:50
$ws.Rows.Item(7).Insert()

# Populate the new row 7 with the new weekly price record
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = Get-Date -Year 2022 -Month 5 -Day 20 -Hour 0 -Minute 0 -Second 0
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112040
$ws.Range("G7").Value = "Cilantro"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 550
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = 575
$ws.Range("N7").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O7").Value = "Provincia de Diguillín"
$ws.Range("P7").Value = 575
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
